# Commit: "changed bharath name to Sreenija"
#   Label.text =”Bharath”        -> Label.text =”Sreenija”
#   Textfield.text = “Kumar”     -> Textfield.text = “Gunnala”
# (only the first "Textfield.text = “Kumar”" paragraph is touched; a later,
# unrelated paragraph with the same "Kumar" text must stay untouched)
#
# A plain Find & Replace keeps the whole paragraph text in a single run.
# The target edit instead leaves the literal/quote text in its own run(s)
# and puts the replaced name in its own run in between - exactly what real
# Word does when you select just the name and retype it. We reproduce that
# here by locating the exact sub-range of the old name via Find (Find
# collapses the range to the match), overwriting its text, and nudging a
# character-formatting property on just that sub-range so it keeps its own
# run instead of being re-merged with its neighbours.

$d = $word.ActiveDocument

function Replace-FirstNameInDoc($oldName, $newName) {
    $target = $d.Content
    $target.Find.Execute($oldName, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
    $target.Text = $newName
    # Touch and revert a formatting property so this run stays split from
    # the surrounding literal-text runs instead of re-merging with them.
    $target.Bold = 1
    $target.Bold = 0
}

Replace-FirstNameInDoc "Bharath" "Sreenija"
Replace-FirstNameInDoc "Kumar" "Gunnala"
